$d = $word.ActiveDocument

$replacements = @(
    @{old = "2025-05-23 Friday"; new = "2025-05-24 Saturday"},
    @{old = "48÷5=9, 3"; new = "81÷9=9, 0"},
    @{old = "14÷2=7, 0"; new = "61÷2=30, 1"},
    @{old = "23÷4=5, 3"; new = "92÷2=46, 0"},
    @{old = "51÷3=17, 0"; new = "55÷8=6, 7"},
    @{old = "72÷4=18, 0"; new = "30÷8=3, 6"},
    @{old = "14÷7=2, 0"; new = "79÷9=8, 7"},
    @{old = "76÷8=9, 4"; new = "74÷8=9, 2"},
    @{old = "50÷4=12, 2"; new = "14÷3=4, 2"},
    @{old = "11÷2=5, 1"; new = "73÷3=24, 1"},
    @{old = "25÷9=2, 7"; new = "62÷2=31, 0"},
    @{old = "60÷7=8, 4"; new = "95÷3=31, 2"},
    @{old = "96÷2=48, 0"; new = "37÷6=6, 1"},
    @{old = "99÷6=16, 3"; new = "86÷4=21, 2"},
    @{old = "26÷9=2, 8"; new = "58÷5=11, 3"},
    @{old = "19÷8=2, 3"; new = "48÷2=24, 0"},
    @{old = "42÷8=5, 2"; new = "76÷7=10, 6"},
    @{old = "57÷2=28, 1"; new = "80÷8=10, 0"},
    @{old = "56÷3=18, 2"; new = "31÷2=15, 1"},
    @{old = "57÷4=14, 1"; new = "59÷7=8, 3"},
    @{old = "33÷9=3, 6"; new = "81÷5=16, 1"},
    @{old = "60÷2=30, 0"; new = "77÷8=9, 5"},
    @{old = "68÷8=8, 4"; new = "64÷8=8, 0"},
    @{old = "93÷4=23, 1"; new = "33÷3=11, 0"},
    @{old = "55÷5=11, 0"; new = "20÷6=3, 2"},
    @{old = "67÷9=7, 4"; new = "98÷7=14, 0"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
